$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows 28 ("SC 92") and 26 ("RM 232").
# Deleting the higher row index first keeps the lower row's index stable.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Column E (missing-data) swaps within rows 2-18 block
$ws.Range("E5").Value = ""
$ws.Range("E8").Value = -6.6
$ws.Range("E12").Value = ""
$ws.Range("E14").Value = -5.4
$ws.Range("E18").Value = ""

# Adjustments in the shifted-up tail rows (26, 27, 33 after the deletions)
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("D33").Value = -14.1
